$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting existing rows 165+ down by one.
$ws.Rows("165:165").Insert()

# Populate the new row 165 with the "Plaude euge Theotocos" entry.
$ws.Range("A165").Value = "plaude_euge_theotocos"
$ws.Range("B165").Value = "Plaude euge Theotocos"
$ws.Range("C165").Value = "Petrus Wilhelmi Grudencz"
$ws.Range("D165").Value = "STB, SBB, ATB, ABB"
$ws.Range("F165").Value = "Past Extoling, Theotokos"
$ws.Range("H165").Value = "Latin"
$ws.Range("I165").Value = "15th century"
$ws.Range("J165").Value = "acrostic"

# Restore a view/selection state close to the author's final saved state
# (the frozen top pane keeps B1 selected, the scrolling bottom pane ends
# up focused on the newly inserted entry's Tags cell).
$win = $excel.ActiveWindow
$win.Panes.Item(1).Activate() | Out-Null
$ws.Range("B1").Select() | Out-Null
$win.Panes.Item(2).Activate() | Out-Null
$ws.Range("J167").Select() | Out-Null
